$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H6").Value = 2842564.2
$ws_ALC.Range("I6").Value = 7579104
$ws_ALC.Range("J6").Value = 640.4
$ws_ALC.Range("K6").Value = 22737312
$ws_ALC.Range("L6").Value = 1921.2
$ws_ALC.Range("M6").Value = -22737200
$ws_ALC.Range("N6").Value = -2145.2

$ws_ALC.Range("H98").Value = 1565.3636
$ws_ALC.Range("I98").Value = 1615.2333
$ws_ALC.Range("J98").Value = 1066.6666
$ws_ALC.Range("K98").Value = 1615.2333
$ws_ALC.Range("L98").Value = 1066.6666
$ws_ALC.Range("M98").Value = -117.2333000000001
$ws_ALC.Range("N98").Value = -4062.6666

$ws_ALC.Range("H112").Value = 3154.1667
$ws_ALC.Range("I112").Value = 1640
$ws_ALC.Range("J112").Value = 3552.6316
$ws_ALC.Range("K112").Value = 4920
$ws_ALC.Range("L112").Value = 10657.8948
$ws_ALC.Range("M112").Value = -3812
$ws_ALC.Range("N112").Value = -12873.8948

$ws_ALC.Range("H116").Value = 4461.8335
$ws_ALC.Range("I116").Value = 4523.0713
$ws_ALC.Range("J116").Value = 4247.5
$ws_ALC.Range("K116").Value = 4523.0713
$ws_ALC.Range("L116").Value = 4247.5
$ws_ALC.Range("M116").Value = -1081.0713
$ws_ALC.Range("N116").Value = -11131.5

$ws_ALC.Range("H122").Value = 1565.3636
$ws_ALC.Range("I122").Value = 1615.2333
$ws_ALC.Range("J122").Value = 1066.6666
$ws_ALC.Range("K122").Value = 4845.699900000001
$ws_ALC.Range("L122").Value = 3199.9998
$ws_ALC.Range("M122").Value = -2395.699900000001
$ws_ALC.Range("N122").Value = -8099.9998

$ws_ALC.Range("H129").Value = 639554.6
$ws_ALC.Range("I129").Value = 437.625
$ws_ALC.Range("J129").Value = 741813.3
$ws_ALC.Range("K129").Value = 1312.875
$ws_ALC.Range("L129").Value = 2225439.9
$ws_ALC.Range("M129").Value = 3687.125
$ws_ALC.Range("N129").Value = -2235439.9

$ws_ALC.Range("H132").Value = 2001749.8
$ws_ALC.Range("I132").Value = 1604.3096
$ws_ALC.Range("J132").Value = 12502513
$ws_ALC.Range("K132").Value = 4812.9288
$ws_ALC.Range("L132").Value = 37507539
$ws_ALC.Range("M132").Value = -2282.9288
$ws_ALC.Range("N132").Value = -37512599

$ws_ALC.Range("H135").Value = 1385.84
$ws_ALC.Range("I135").Value = 1421.238
$ws_ALC.Range("J135").Value = 1200
$ws_ALC.Range("K135").Value = 12791.142
$ws_ALC.Range("L135").Value = 10800
$ws_ALC.Range("M135").Value = -10256.142
$ws_ALC.Range("N135").Value = -15870

$ws_ALC.Range("H137").Value = 1008.88464
$ws_ALC.Range("I137").Value = 911.2381
$ws_ALC.Range("J137").Value = 1419
$ws_ALC.Range("K137").Value = 2733.7143
$ws_ALC.Range("L137").Value = 4257
$ws_ALC.Range("M137").Value = -183.7143000000001
$ws_ALC.Range("N137").Value = -9357

$ws_ALC.Range("H141").Value = 742
$ws_ALC.Range("I141").Value = 684.0909
$ws_ALC.Range("J141").Value = 901.25
$ws_ALC.Range("K141").Value = 2052.2727
$ws_ALC.Range("L141").Value = 2703.75
$ws_ALC.Range("M141").Value = 3127.7273
$ws_ALC.Range("N141").Value = -13063.75

$ws_ARM.Range("H3").Value = 200
$ws_ARM.Range("I3").Value = 200
$ws_ARM.Range("J3").Value = 0
$ws_ARM.Range("K3").Value = 200
$ws_ARM.Range("L3").Value = 0
$ws_ARM.Range("M3").Value = -85
$ws_ARM.Range("N3").ClearContents()

$ws_ARM.Range("H32").Value = 17785.59
$ws_ARM.Range("I32").Value = 13015.414
$ws_ARM.Range("J32").Value = 110009
$ws_ARM.Range("K32").Value = 13015.414
$ws_ARM.Range("L32").Value = 110009
$ws_ARM.Range("M32").Value = -12728.414
$ws_ARM.Range("N32").Value = -110583

$ws_ARM.Range("H61").Value = 1885.0625
$ws_ARM.Range("I61").Value = 1561.5
$ws_ARM.Range("J61").Value = 4150
$ws_ARM.Range("K61").Value = 1561.5
$ws_ARM.Range("L61").Value = 4150
$ws_ARM.Range("M61").Value = -1349.5
$ws_ARM.Range("N61").Value = -4574

$ws_ARM.Range("H132").Value = 2550.5862
$ws_ARM.Range("I132").Value = 1947.7778
$ws_ARM.Range("J132").Value = 3537
$ws_ARM.Range("K132").Value = 5843.3334
$ws_ARM.Range("L132").Value = 10611
$ws_ARM.Range("M132").Value = -3313.3334
$ws_ARM.Range("N132").Value = -15671

$ws_ARM.Range("H136").Value = 1885.0625
$ws_ARM.Range("I136").Value = 1561.5
$ws_ARM.Range("J136").Value = 4150
$ws_ARM.Range("K136").Value = 4684.5
$ws_ARM.Range("L136").Value = 12450
$ws_ARM.Range("M136").Value = -2134.5
$ws_ARM.Range("N136").Value = -17550

$ws_BSM.Range("H11").Value = 100005
$ws_BSM.Range("I11").Value = 0
$ws_BSM.Range("J11").Value = 100005
$ws_BSM.Range("K11").Value = 0
$ws_BSM.Range("L11").Value = 100005
$ws_BSM.Range("N11").Value = -100285

$ws_BSM.Range("H134").Value = 20805.076
$ws_BSM.Range("I134").Value = 1644.75
$ws_BSM.Range("J134").Value = 114477.78
$ws_BSM.Range("K134").Value = 4934.25
$ws_BSM.Range("L134").Value = 343433.34
$ws_BSM.Range("M134").Value = -2399.25
$ws_BSM.Range("N134").Value = -348503.34

$ws_CRP.Range("H3").Value = 3980
$ws_CRP.Range("I3").Value = 0
$ws_CRP.Range("J3").Value = 3980
$ws_CRP.Range("K3").Value = 0
$ws_CRP.Range("L3").Value = 3980
$ws_CRP.Range("N3").Value = -4206

$ws_CRP.Range("H31").Value = 4066.3044
$ws_CRP.Range("I31").Value = 3969.1428
$ws_CRP.Range("J31").Value = 4217.4443
$ws_CRP.Range("K31").Value = 3969.1428
$ws_CRP.Range("L31").Value = 4217.4443
$ws_CRP.Range("M31").Value = -3674.1428
$ws_CRP.Range("N31").Value = -4807.4443

$ws_CRP.Range("H34").Value = 4066.3044
$ws_CRP.Range("I34").Value = 3969.1428
$ws_CRP.Range("J34").Value = 4217.4443
$ws_CRP.Range("K34").Value = 3969.1428
$ws_CRP.Range("L34").Value = 4217.4443
$ws_CRP.Range("M34").Value = -3767.1428
$ws_CRP.Range("N34").Value = -4621.4443

$ws_CRP.Range("H132").Value = 1372.4849
$ws_CRP.Range("I132").Value = 1014.76
$ws_CRP.Range("J132").Value = 2490.375
$ws_CRP.Range("K132").Value = 3044.28
$ws_CRP.Range("L132").Value = 7471.125
$ws_CRP.Range("M132").Value = -514.2799999999997
$ws_CRP.Range("N132").Value = -12531.125

$ws_CRP.Range("H141").Value = 56800.4
$ws_CRP.Range("I141").Value = 23000
$ws_CRP.Range("J141").Value = 60556
$ws_CRP.Range("K141").Value = 23000
$ws_CRP.Range("L141").Value = 60556
$ws_CRP.Range("M141").Value = -17820
$ws_CRP.Range("N141").Value = -70916

$ws_CUL.Range("H6").Value = 105.666664
$ws_CUL.Range("I6").Value = 91.77778000000001
$ws_CUL.Range("J6").Value = 147.33333
$ws_CUL.Range("K6").Value = 275.33334
$ws_CUL.Range("L6").Value = 441.99999
$ws_CUL.Range("M6").Value = -162.33334
$ws_CUL.Range("N6").Value = -667.99999

$ws_CUL.Range("H43").Value = 0
$ws_CUL.Range("I43").Value = 0
$ws_CUL.Range("J43").Value = 0
$ws_CUL.Range("K43").Value = 0
$ws_CUL.Range("L43").Value = 0
$ws_CUL.Range("N43").ClearContents()

$ws_CUL.Range("H113").Value = 742.6111
$ws_CUL.Range("I113").Value = 1305.25
$ws_CUL.Range("J113").Value = 581.8570999999999
$ws_CUL.Range("K113").Value = 3915.75
$ws_CUL.Range("L113").Value = 1745.5713
$ws_CUL.Range("M113").Value = -1745.75
$ws_CUL.Range("N113").Value = -6085.5713

$ws_CUL.Range("H122").Value = 22953.674
$ws_CUL.Range("I122").Value = 593
$ws_CUL.Range("J122").Value = 26967.129
$ws_CUL.Range("K122").Value = 5337
$ws_CUL.Range("L122").Value = 242704.161
$ws_CUL.Range("M122").Value = -2887
$ws_CUL.Range("N122").Value = -247604.161

$ws_GSM.Range("H7").Value = 0
$ws_GSM.Range("I7").Value = 0
$ws_GSM.Range("J7").Value = 0
$ws_GSM.Range("K7").Value = 0
$ws_GSM.Range("L7").Value = 0
$ws_GSM.Range("M7").ClearContents()

$ws_GSM.Range("H8").Value = 0
$ws_GSM.Range("I8").Value = 0
$ws_GSM.Range("J8").Value = 0
$ws_GSM.Range("K8").Value = 0
$ws_GSM.Range("L8").Value = 0
$ws_GSM.Range("M8").ClearContents()

$ws_WVR.Range("H132").Value = 9680.857
$ws_WVR.Range("I132").Value = 3176.8
$ws_WVR.Range("J132").Value = 13294.223
$ws_WVR.Range("K132").Value = 9530.400000000001
$ws_WVR.Range("L132").Value = 39882.669
$ws_WVR.Range("M132").Value = -7000.400000000001
$ws_WVR.Range("N132").Value = -44942.669

$ws_WVR.Range("H133").Value = 44091.3
$ws_WVR.Range("I133").Value = 0
$ws_WVR.Range("J133").Value = 44091.3
$ws_WVR.Range("K133").Value = 0
$ws_WVR.Range("L133").Value = 44091.3
$ws_WVR.Range("N133").Value = -54211.3

$ws_WVR.Range("H136").Value = 4752.2964
$ws_WVR.Range("I136").Value = 853.41174
$ws_WVR.Range("J136").Value = 11380.4
$ws_WVR.Range("K136").Value = 2560.23522
$ws_WVR.Range("L136").Value = 34141.2
$ws_WVR.Range("M136").Value = -10.23522000000003
$ws_WVR.Range("N136").Value = -39241.2
